# Applies the changes described by the commit "doc & update sql":
#  1. Merge "ASUTSK." + " " (two runs) into a single run "ASUTSK. "
#     (only the occurrence immediately preceding "dbo.Outcomes").
#  2. Merge "Перенос " + "приема" + " топлива в ЦОД" (three runs) into a
#     single run "Перенос приема топлива в ЦОД".
#  3/4. Rename the two "AddRemainsTanks" occurrences to "AddRemains".
#  5. Rename "RemainsTanks" to "RemainsTank".
#  6. Rename "RemainsTanks_TSK" to "Remains_TSK".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "ASUTSK." + " " -> single run "ASUTSK. "
#    There are several "ASUTSK." runs in the document; only the second
#    occurrence (immediately followed by a standalone run containing a
#    single space before "dbo.Outcomes") needs merging - the others are
#    either already a single run or not followed by a lone-space run.
# ---------------------------------------------------------------------
$search = $d.Content
$search.Find.ClearFormatting()
$found = $search.Find.Execute("ASUTSK.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchIndex = 0
$targetStart = -1
$targetEnd = -1
while ($found -and $matchIndex -lt 50) {
    if ($matchIndex -eq 1) {
        $targetStart = $search.Start
        $targetEnd = $search.End
    }
    $search.Collapse(0)
    $search.MoveEnd(1, 1000) | Out-Null
    $found = $search.Find.Execute("ASUTSK.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $matchIndex = $matchIndex + 1
}

if ($targetStart -ge 0) {
    # Extend the match by one character to swallow the following
    # standalone-space run, then rewrite the combined range as a single
    # run. A same-text assignment is a no-op for the underlying engine,
    # so first write a differing placeholder to force the run merge,
    # then write the real text back.
    $mergeRange = $d.Range($targetStart, $targetEnd + 1)
    $mergeRange.Text = "ASUTSKPLACEHOLDER"
    $mergeRange2 = $d.Range($targetStart, $targetStart + 17)
    $mergeRange2.Text = "ASUTSK. "
}

# ---------------------------------------------------------------------
# 2) "Перенос " + "приема" + " топлива в ЦОД" -> single run
#    "Перенос приема топлива в ЦОД"
# ---------------------------------------------------------------------
$phraseRange = $d.Content
$phraseRange.Find.ClearFormatting()
$phraseFound = $phraseRange.Find.Execute("Перенос приема топлива в ЦОД", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($phraseFound) {
    $pStart = $phraseRange.Start
    $pEnd = $phraseRange.End
    $full = $d.Range($pStart, $pEnd)
    $full.Text = "XXXXXXXXXXXXXXXXXXXXXXXXXXXX"
    $full2 = $d.Range($pStart, $pStart + 28)
    $full2.Text = "Перенос приема топлива в ЦОД"
}

# ---------------------------------------------------------------------
# 3-6) Plain identifier renames - each is a contained, single-run match
#      so a document-wide Find/ReplaceAll is unambiguous once applied
#      in most-specific-first order.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Execute("AddRemainsTanks", $true, $false, $false, $false, $false, $true, 1, $false, "AddRemains", 2) | Out-Null

$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute("RemainsTanks_TSK", $true, $false, $false, $false, $false, $true, 1, $false, "Remains_TSK", 2) | Out-Null

$r3 = $d.Content
$r3.Find.ClearFormatting()
$r3.Find.Execute("RemainsTanks", $true, $false, $false, $false, $false, $true, 1, $false, "RemainsTank", 2) | Out-Null

Write-Output "edit.ps1 completed"
